$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.898150666666666
$ws.Range("H2").Value = 26.694452
$ws.Range("I2").Value = 0.3765197173862137
$ws.Range("J2").Value = 0.3765197173862137
$ws.Range("M2").Value = 0.3284223333333333
$ws.Range("N2").Value = 0.985267
$ws.Range("O2").Value = 0.04541528350839906
$ws.Range("P2").Value = 0.04541528350839906
$ws.Range("Q2").Value = 2.922351404298222
$ws.Range("R2").Value = 26.301162638684
$ws.Range("S2").Value = 0.01709974971159718
$ws.Range("T2").Value = 0.01709974971159718

$ws.Range("G3").Value = 8.898150666666666
$ws.Range("H3").Value = 26.694452
$ws.Range("I3").Value = 0.3765197173862137
$ws.Range("J3").Value = 0.3765197173862137
$ws.Range("M3").Value = 4.062688333333334
$ws.Range("N3").Value = 12.188065
$ws.Range("O3").Value = 0.5618014481290817
$ws.Range("P3").Value = 0.5618014481290816
$ws.Range("Q3").Value = 36.15041290170889
$ws.Range("R3").Value = 325.3537161153801
$ws.Range("S3").Value = 0.2115293224767275
$ws.Range("T3").Value = 0.2115293224767274

$ws.Range("G4").Value = 8.898150666666666
$ws.Range("H4").Value = 26.694452
$ws.Range("I4").Value = 0.3765197173862137
$ws.Range("J4").Value = 0.3765197173862137
$ws.Range("M4").Value = 2.840427
$ws.Range("N4").Value = 8.521281
$ws.Range("O4").Value = 0.3927832683625193
$ws.Range("P4").Value = 0.3927832683625193
$ws.Range("Q4").Value = 25.274547403668
$ws.Range("R4").Value = 227.470926633012
$ws.Range("S4").Value = 0.1478906451978891
$ws.Range("T4").Value = 0.1478906451978891

$ws.Range("I5").Value = 0.1415167724465014
$ws.Range("J5").Value = 0.1415167724465015
$ws.Range("M5").Value = 0.3284223333333333
$ws.Range("N5").Value = 0.985267
$ws.Range("O5").Value = 0.04541528350839906
$ws.Range("P5").Value = 0.04541528350839906
$ws.Range("Q5").Value = 1.098380030564444
$ws.Range("R5").Value = 9.885420275080001
$ws.Range("S5").Value = 0.006427024341851459
$ws.Range("T5").Value = 0.00642702434185146

$ws.Range("I6").Value = 0.1415167724465014
$ws.Range("J6").Value = 0.1415167724465015
$ws.Range("M6").Value = 4.062688333333334
$ws.Range("N6").Value = 12.188065
$ws.Range("O6").Value = 0.5618014481290817
$ws.Range("P6").Value = 0.5618014481290816
$ws.Range("S6").Value = 0.07950432769499824
$ws.Range("T6").Value = 0.07950432769499824

$ws.Range("I7").Value = 0.1415167724465014
$ws.Range("J7").Value = 0.1415167724465015
$ws.Range("M7").Value = 2.840427
$ws.Range("N7").Value = 8.521281
$ws.Range("O7").Value = 0.3927832683625193
$ws.Range("P7").Value = 0.3927832683625193
$ws.Range("Q7").Value = 9.49956193116
$ws.Range("R7").Value = 85.49605738044
$ws.Range("S7").Value = 0.05558542040965175
$ws.Range("T7").Value = 0.05558542040965175

$ws.Range("G8").Value = 11.39006466666667
$ws.Range("H8").Value = 34.170194
$ws.Range("I8").Value = 0.4819635101672848
$ws.Range("J8").Value = 0.4819635101672848
$ws.Range("M8").Value = 0.3284223333333333
$ws.Range("N8").Value = 0.985267
$ws.Range("O8").Value = 0.04541528350839906
$ws.Range("P8").Value = 0.04541528350839906
$ws.Range("Q8").Value = 3.740751614644222
$ws.Range("R8").Value = 33.666764531798
$ws.Range("S8").Value = 0.02188850945495041
$ws.Range("T8").Value = 0.02188850945495041

$ws.Range("G9").Value = 11.39006466666667
$ws.Range("H9").Value = 34.170194
$ws.Range("I9").Value = 0.4819635101672848
$ws.Range("J9").Value = 0.4819635101672848
$ws.Range("M9").Value = 4.062688333333334
$ws.Range("N9").Value = 12.188065
$ws.Range("O9").Value = 0.5618014481290817
$ws.Range("P9").Value = 0.5618014481290816
$ws.Range("Q9").Value = 46.2742828371789
$ws.Range("R9").Value = 416.4685455346101
$ws.Range("S9").Value = 0.270767797957356
$ws.Range("T9").Value = 0.2707677979573559

$ws.Range("G10").Value = 11.39006466666667
$ws.Range("H10").Value = 34.170194
$ws.Range("I10").Value = 0.4819635101672848
$ws.Range("J10").Value = 0.4819635101672848
$ws.Range("M10").Value = 2.840427
$ws.Range("N10").Value = 8.521281
$ws.Range("O10").Value = 0.3927832683625193
$ws.Range("P10").Value = 0.3927832683625193
$ws.Range("Q10").Value = 32.352647210946
$ws.Range("R10").Value = 291.173824898514
$ws.Range("S10").Value = 0.1893072027549784
$ws.Range("T10").Value = 0.1893072027549784
